$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new (blank) column before column N ("Late" / Repayment schedule sheet),
# shifting the existing N:R columns one position to the right (-> O:S).
$ws.Columns("N").Insert()

# New column N should take on the same width as column M ("Principal" column,
# width 10.7109375 chars) rather than the sheet default width.
$ws.Columns("N").ColumnWidth = 9.877604166666666

# Move the active selection to match the post-edit state.
[void]$ws.Range("S7").Select()
